# -----------------------------------------------------------------------
# Edit 1: "Mt. Ulap" -> split into two runs ("Mt. " / "Ulap") with
# spell-check proofErr markers bracketing "Ulap" (spellStart / spellEnd),
# matching Word's normal behaviour when the spell checker flags a word.
# -----------------------------------------------------------------------
$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("Ulap", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found) {
    $oldStart = $rng.Start
    $oldEnd   = $rng.End

    # Insert the replacement run (wrapped in proofErr markers) right after
    # the found "Ulap" text; the surrounding run "Mt. " is left untouched.
    $xmlUlap = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ulap</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xmlUlap)

    # Now remove the original (un-marked) "Ulap" text that is still sitting
    # just before the text we inserted.
    $rngOld = $d.Range($oldStart, $oldEnd)
    $rngOld.Text = ""
}

# -----------------------------------------------------------------------
# Edit 2: after the "Resorts" paragraph, add a blank paragraph, a new
# "Iloilo" paragraph, and another blank paragraph (the pre-existing blank
# paragraph right before the sectPr is left as-is).
# -----------------------------------------------------------------------
$last = $d.Paragraphs.Last
$tailRng = $last.Range

$xmlIloilo = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:r><w:t>Iloilo</w:t></w:r></w:p><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$tailRng.InsertXML($xmlIloilo)

Write-Output "Done"
